# Update the "2019" servo-angle tracking sheet with the newly measured
# values for the Spine / Omoplate servos and mark the Arm-Left rotation
# rows (Rotate / Shoulder / Omoplate) as "done" with the green fill used
# for every other completed body part.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")

# --- Spine row (row 35): new min/max angle + measured value -------------
$ws.Range("C35").Value = 80
$ws.Range("D35").Value = 125
$ws.Range("F35").Value = 85

# --- Omoplate row (row 36): new min/max angle + measured value ----------
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 0
$ws.Range("F36").Value = 40

# --- Mark this servo group ("Rotate", "Shoulder", "Omoplate" labels in
#     column A, rows 34-36) as configured, matching the green highlight
#     already used on the other completed rows (e.g. A6:A7) ------------
$ws.Range("A34:A36").Interior.Color = 5296274

# --- Update the sheet's view/selection state to reflect where the user
#     left off working (best effort - scroll so row 6 is at the top and
#     leave the selection on the last edited cell) -----------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A36").Select()

Write-Output "Updated 2019 servo sheet: Spine/Omoplate angles + Arm-Left rotation highlight"
